# The underlying dataset gained one additional weekly observation.
# A new row is inserted at sheet row 737 (pushing the existing rows
# 737-804 down to 738-805) and populated with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 737; Excel shifts 737:804 down to 738:805
# and the dimension/used range grows to A1:R805 automatically.
$ws.Rows(737).Insert()

# Populate the newly inserted row 737 with the new observation.
$ws.Cells.Item(737, 1).Value  = 6
$ws.Cells.Item(737, 2).Value  = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(737, 3).Value  = "Metropolitana"
$ws.Cells.Item(737, 4).Value  = 45166
$ws.Cells.Item(737, 5).Value  = 13
$ws.Cells.Item(737, 6).Value  = 100112044
$ws.Cells.Item(737, 7).Value  = "Perejil"
$ws.Cells.Item(737, 8).Value  = "Sin especificar"
$ws.Cells.Item(737, 9).Value  = "Primera"
$ws.Cells.Item(737, 10).Value = 270
$ws.Cells.Item(737, 11).Value = 9000
$ws.Cells.Item(737, 12).Value = 10000
$ws.Cells.Item(737, 13).Value = 9667
$ws.Cells.Item(737, 14).Value = "`$/docena de atados"
$ws.Cells.Item(737, 15).Value = "Región Metropolitana"
$ws.Cells.Item(737, 16).Value = 3222
$ws.Cells.Item(737, 17).Value = 3
$ws.Cells.Item(737, 18).Value = "Hortaliza"
